$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column G: the ExpectedFilenames list now has 3 rows per category
# (StandardExcelReport-...-2023_, ExcelReport-Takeda-...-, WordReport-...-)
# instead of 2 (ExcelReport-..., WordReport-...), so the column grows from
# G2:G10 to G2:G13.
$ws.Range("G2").Value = "StandardExcelReport-Takeda - MM Maintenance-Clinical-2023_"
$ws.Range("G3").Value = "ExcelReport-Takeda-MM Maintenance-Clinical-"
$ws.Range("G4").Value = "WordReport-Takeda - MM Maintenance-Clinical-"
$ws.Range("G5").Value = "StandardExcelReport-Takeda - MM Maintenance-Economic-2023_"
$ws.Range("G6").Value = "ExcelReport-Takeda-MM Maintenance-Economic-"
$ws.Range("G7").Value = "WordReport-Takeda - MM Maintenance-Economic-"
$ws.Range("G8").Value = "StandardExcelReport-Takeda - MM Maintenance-Quality of Life-2023_"
$ws.Range("G9").Value = "ExcelReport-Takeda-MM Maintenance-Quality of Life-"
$ws.Range("G10").Value = "WordReport-Takeda - MM Maintenance-Quality of Life-"
$ws.Range("G11").Value = "StandardExcelReport-Takeda - MM Maintenance-Real-world Evidence-2023_"
$ws.Range("G12").Value = "ExcelReport-Takeda-MM Maintenance-Real-world Evidence-"
$ws.Range("G13").Value = "WordReport-Takeda - MM Maintenance-Real-world Evidence-"

# Clear the bestFit/center-align styling that used to be applied to A2:A5,
# D3:D5, E2:E5, F2:F5 - the data rows now use the plain default style.
$ws.Range("A2:A5").Style = "Normal"
$ws.Range("D3:D5").Style = "Normal"
$ws.Range("E2:E5").Style = "Normal"
$ws.Range("F2:F5").Style = "Normal"

# Selection/view now anchored further right, covering the newly extended G column.
$ws.Range("G2:G13").Select()
$ws.Application.ActiveWindow.ScrollColumn = 6
